# chore: update Sheets via scheduled runner
# Refreshes cached market-board pricing/profit figures (currentAveragePrice*,
# LevePrice*, LeveProfit*) for a batch of leve rows across all job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 2351.9092
$ws.Range("I8").Value = 2351.9092
$ws.Range("K8").Value = 7055.7276
$ws.Range("M8").Value = -6916.7276

$ws.Range("H19").Value = 1800.8
$ws.Range("I19").Value = 2527.7144
$ws.Range("J19").Value = 1164.75
$ws.Range("K19").Value = 2527.7144
$ws.Range("L19").Value = 1164.75
$ws.Range("M19").Value = -2352.7144
$ws.Range("N19").Value = -1514.75

$ws.Range("H28").Value = 1607.8125
$ws.Range("I28").Value = 2049.375
$ws.Range("J28").Value = 1166.25
$ws.Range("K28").Value = 2049.375
$ws.Range("L28").Value = 1166.25
$ws.Range("M28").Value = -1564.375
$ws.Range("N28").Value = -2136.25

$ws.Range("H33").Value = 8177.6333
$ws.Range("I33").Value = 1871.6957
$ws.Range("J33").Value = 28897.143
$ws.Range("K33").Value = 1871.6957
$ws.Range("L33").Value = 28897.143
$ws.Range("M33").Value = -1642.6957
$ws.Range("N33").Value = -29355.143

$ws.Range("H62").Value = 1223
$ws.Range("I62").Value = 1165.9166
$ws.Range("J62").Value = 1360
$ws.Range("K62").Value = 1165.9166
$ws.Range("L62").Value = 1360
$ws.Range("M62").Value = -541.9166
$ws.Range("N62").Value = -2608

$ws.Range("H65").Value = 1223
$ws.Range("I65").Value = 1165.9166
$ws.Range("J65").Value = 1360
$ws.Range("K65").Value = 5829.583000000001
$ws.Range("L65").Value = 6800
$ws.Range("M65").Value = -2709.583000000001
$ws.Range("N65").Value = -13040

$ws.Range("H88").Value = 650.5263
$ws.Range("I88").Value = 761.2
$ws.Range("J88").Value = 611
$ws.Range("K88").Value = 761.2
$ws.Range("L88").Value = 611
$ws.Range("M88").Value = -355.2
$ws.Range("N88").Value = -1423

$ws.Range("H91").Value = 650.5263
$ws.Range("I91").Value = 761.2
$ws.Range("J91").Value = 611
$ws.Range("K91").Value = 761.2
$ws.Range("L91").Value = 611
$ws.Range("M91").Value = 642.8
$ws.Range("N91").Value = -3419

$ws.Range("H98").Value = 1541.6428
$ws.Range("I98").Value = 1068.625
$ws.Range("J98").Value = 2172.3333
$ws.Range("K98").Value = 1068.625
$ws.Range("L98").Value = 2172.3333
$ws.Range("M98").Value = 429.375
$ws.Range("N98").Value = -5168.3333

$ws.Range("H107").Value = 884.2069
$ws.Range("I107").Value = 513.17645
$ws.Range("J107").Value = 1409.8334
$ws.Range("K107").Value = 513.17645
$ws.Range("L107").Value = 1409.8334
$ws.Range("M107").Value = 1406.82355
$ws.Range("N107").Value = -5249.8334

$ws.Range("H116").Value = 2559867
$ws.Range("I116").Value = 23811856
$ws.Range("J116").Value = 9628.200000000001
$ws.Range("K116").Value = 23811856
$ws.Range("L116").Value = 9628.200000000001
$ws.Range("M116").Value = -23808414
$ws.Range("N116").Value = -16512.2

$ws.Range("H122").Value = 1541.6428
$ws.Range("I122").Value = 1068.625
$ws.Range("J122").Value = 2172.3333
$ws.Range("K122").Value = 3205.875
$ws.Range("L122").Value = 6516.999899999999
$ws.Range("M122").Value = -755.875
$ws.Range("N122").Value = -11416.9999

$ws.Range("H132").Value = 2166311
$ws.Range("I132").Value = 3135839
$ws.Range("J132").Value = 3517.923
$ws.Range("K132").Value = 9407517
$ws.Range("L132").Value = 10553.769
$ws.Range("M132").Value = -9404987
$ws.Range("N132").Value = -15613.769

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2063.5151
$ws.Range("I61").Value = 1511.6522
$ws.Range("J61").Value = 3332.8
$ws.Range("K61").Value = 1511.6522
$ws.Range("L61").Value = 3332.8
$ws.Range("M61").Value = -1299.6522
$ws.Range("N61").Value = -3756.8

$ws.Range("H63").Value = 1682.25
$ws.Range("I63").Value = 1515.3125
$ws.Range("J63").Value = 2350
$ws.Range("K63").Value = 1515.3125
$ws.Range("L63").Value = 2350
$ws.Range("M63").Value = -829.3125
$ws.Range("N63").Value = -3722

$ws.Range("H66").Value = 1682.25
$ws.Range("I66").Value = 1515.3125
$ws.Range("J66").Value = 2350
$ws.Range("K66").Value = 7576.5625
$ws.Range("L66").Value = 11750
$ws.Range("M66").Value = -4144.5625
$ws.Range("N66").Value = -18614

$ws.Range("H119").Value = 40174.375
$ws.Range("J119").Value = 40174.375
$ws.Range("L119").Value = 40174.375
$ws.Range("N119").Value = -49850.375

$ws.Range("H136").Value = 2063.5151
$ws.Range("I136").Value = 1511.6522
$ws.Range("J136").Value = 3332.8
$ws.Range("K136").Value = 4534.9566
$ws.Range("L136").Value = 9998.400000000001
$ws.Range("M136").Value = -1984.9566
$ws.Range("N136").Value = -15098.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1815.3636
$ws.Range("I134").Value = 1401.4117
$ws.Range("J134").Value = 3222.8
$ws.Range("K134").Value = 4204.2351
$ws.Range("L134").Value = 9668.400000000001
$ws.Range("M134").Value = -1669.2351
$ws.Range("N134").Value = -14738.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2097.5483
$ws.Range("I58").Value = 958.06665
$ws.Range("J58").Value = 3165.8125
$ws.Range("K58").Value = 958.06665
$ws.Range("L58").Value = 3165.8125
$ws.Range("M58").Value = -755.06665
$ws.Range("N58").Value = -3571.8125

$ws.Range("H62").Value = 168101
$ws.Range("I62").Value = 201320
$ws.Range("J62").Value = 2006
$ws.Range("K62").Value = 201320
$ws.Range("L62").Value = 2006
$ws.Range("M62").Value = -200696
$ws.Range("N62").Value = -3254

$ws.Range("H65").Value = 168101
$ws.Range("I65").Value = 201320
$ws.Range("J65").Value = 2006
$ws.Range("K65").Value = 1006600
$ws.Range("L65").Value = 10030
$ws.Range("M65").Value = -1003480
$ws.Range("N65").Value = -16270

$ws.Range("H122").Value = 3925.6511
$ws.Range("I122").Value = 7036.1177
$ws.Range("J122").Value = 1891.8846
$ws.Range("K122").Value = 21108.3531
$ws.Range("L122").Value = 5675.6538
$ws.Range("M122").Value = -18658.3531
$ws.Range("N122").Value = -10575.6538

$ws.Range("H132").Value = 2040.2084
$ws.Range("I132").Value = 1429.6666
$ws.Range("K132").Value = 4288.9998
$ws.Range("M132").Value = -1758.9998

$ws.Range("H134").Value = 1382.2941
$ws.Range("I134").Value = 984.6429000000001
$ws.Range("K134").Value = 2953.9287
$ws.Range("M134").Value = -418.9287000000004

$ws.Range("H136").Value = 2097.5483
$ws.Range("I136").Value = 958.06665
$ws.Range("J136").Value = 3165.8125
$ws.Range("K136").Value = 2874.19995
$ws.Range("L136").Value = 9497.4375
$ws.Range("M136").Value = -324.1999500000002
$ws.Range("N136").Value = -14597.4375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 7578.643
$ws.Range("I68").Value = 25150.5
$ws.Range("K68").Value = 75451.5
$ws.Range("M68").Value = -74640.5

$ws.Range("H71").Value = 7578.643
$ws.Range("I71").Value = 25150.5
$ws.Range("K71").Value = 226354.5
$ws.Range("M71").Value = -222298.5

$ws.Range("H131").Value = 1523.9799
$ws.Range("I131").Value = 293.46667
$ws.Range("J131").Value = 1743.7142
$ws.Range("K131").Value = 880.4000100000001
$ws.Range("L131").Value = 5231.142599999999
$ws.Range("M131").Value = 4159.59999
$ws.Range("N131").Value = -15311.1426

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H69").Value = 53800
$ws.Range("J69").Value = 53800
$ws.Range("L69").Value = 53800
$ws.Range("N69").Value = -55298

$ws.Range("H70").Value = 4795.8
$ws.Range("I70").Value = 5254
$ws.Range("J70").Value = 4490.3335
$ws.Range("K70").Value = 5254
$ws.Range("L70").Value = 4490.3335
$ws.Range("M70").Value = -4984
$ws.Range("N70").Value = -5030.3335

$ws.Range("H72").Value = 53800
$ws.Range("J72").Value = 53800
$ws.Range("L72").Value = 161400
$ws.Range("N72").Value = -168888

$ws.Range("H73").Value = 4795.8
$ws.Range("I73").Value = 5254
$ws.Range("J73").Value = 4490.3335
$ws.Range("K73").Value = 5254
$ws.Range("L73").Value = 4490.3335
$ws.Range("M73").Value = -4318
$ws.Range("N73").Value = -6362.3335

$ws.Range("H102").Value = 2446
$ws.Range("I102").Value = 1548.5454
$ws.Range("K102").Value = 1548.5454
$ws.Range("M102").Value = 73.45460000000003

$ws.Range("H123").Value = 20199.166
$ws.Range("J123").Value = 20199.166
$ws.Range("L123").Value = 20199.166
$ws.Range("N123").Value = -25099.166

$ws.Range("H132").Value = 8221.105
$ws.Range("I132").Value = 14090.5
$ws.Range("J132").Value = 3952.4546
$ws.Range("K132").Value = 42271.5
$ws.Range("L132").Value = 11857.3638
$ws.Range("M132").Value = -39741.5
$ws.Range("N132").Value = -16917.3638

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 16317.556
$ws.Range("I122").Value = 41002.668
$ws.Range("J122").Value = 3975
$ws.Range("K122").Value = 123008.004
$ws.Range("L122").Value = 11925
$ws.Range("M122").Value = -120558.004
$ws.Range("N122").Value = -16825

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 11666.667
$ws.Range("J31").Value = 11666.667
$ws.Range("L31").Value = 11666.667
$ws.Range("N31").Value = -12362.667

$ws.Range("H113").Value = 473.84848
$ws.Range("I113").Value = 157.05263
$ws.Range("K113").Value = 471.15789
$ws.Range("M113").Value = 1698.84211
